$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = 6.2184908458334
$ws.Range("E4").Value = 2.284984018611098
$ws.Range("F4").Value = 0.212494666666656
$ws.Range("G4").Value = 2.497478685277773
